# The commit inserts one new data row into the "Hortaliza, Vega Monumental
# Concepción - Repollo" sheet at row 530, pushing the existing rows 530-609
# down to 531-610 (dimension grows from A1:R609 to A1:R610).
#
# Reproduce that with a real row insert (so styles/formatting on the
# existing rows shift down exactly like Excel's Insert does), then fill
# in the values for the brand-new row 530.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 530; Excel pushes old row 530 (and everything
# below it) down by one, which is exactly the behaviour the diff shows.
$ws.Rows.Item(530).Insert()

# Populate the newly inserted row 530 with the new record.
$ws.Range("A530").Value = 11
$ws.Range("B530").Value = "Vega Monumental Concepción"
$ws.Range("C530").Value = "Bíobío"
$ws.Range("D530").Value = 45127
$ws.Range("E530").Value = 8
$ws.Range("F530").Value = 100112006
$ws.Range("G530").Value = "Repollo"
$ws.Range("H530").Value = "Copenhague"
$ws.Range("I530").Value = "Primera"
$ws.Range("J530").Value = 1500
$ws.Range("K530").Value = 800
$ws.Range("L530").Value = 1000
$ws.Range("M530").Value = 933
$ws.Range("N530").Value = "$/unidad"
$ws.Range("O530").Value = "Región Metropolitana"
$ws.Range("P530").Value = 933
$ws.Range("Q530").Value = 1
$ws.Range("R530").Value = "Hortaliza"
